$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "configuration path shall default..." requirement row (row 21).
$ws.Rows.Item(21).Delete()

# Make room for two new requirement rows describing the config/spec file
# defaulting and directory-search behaviour.
$ws.Rows.Item(21).Resize(2).Insert()

# Row 19 unchanged in wording, but now carries a numeric "11" in column A.
$ws.Range("A19").Value = 11

# Row 20: reworded from "shall be supported" to "shall be required", also
# gains the numeric "11" marker in column A.
$ws.Range("A20").Value = 11
$ws.Range("C20").Value = "Configuration specification shall be required"

# Row 21 (new): default file naming behaviour.
$ws.Range("B21").Value = "Configuration"
$ws.Range("C21").Value = "providing no file names will default to 'config.ini' and 'config_spec.ini', for the configuration and specification files, respectively"

# Row 22 (new): directory search behaviour.
$ws.Range("B22").Value = "Configuration"
$ws.Range("C22").Value = "if configuration files are not found in the local directory, each higher directory shall be searched for the files until a specific folder is reached"

# Row 23 (previously the "parameter access..." row) gains the "11" marker.
$ws.Range("A23").Value = 11

# Row 24 (previously "A configured parameter may be optional or required")
# gains the "11" marker.
$ws.Range("A24").Value = 11

# Grow the requirements table to account for the extra row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C77"))

# Match the author's final selection.
$ws.Range("C24").Select()
